# 3DES Projetos - Adicionada aula08 - Desavio
# Adds a new week block (columns Q:U) to the FREQ sheet, mirroring the
# existing PROJ/PROJ/PROJ/PDMO/RMST header + weekday-date pattern used
# by the three prior weeks (B:F, G:K, L:P), fills in the new
# attendance (P/F) column for every student, normalizes the existing
# P column away from the lowercase "p"/"f" duplicates onto the
# canonical "P"/"F" strings, and freezes the header row/name column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")

# --- Week 4 header row (row 1): PROJ, PROJ, PROJ, PDMO, RMST ---
$ws.Range("Q1").Value = "PROJ"
$ws.Range("R1").Value = "PROJ"
$ws.Range("S1").Value = "PROJ"
$ws.Range("T1").Value = "PDMO"
$ws.Range("U1").Value = "RMST"

# --- Week 4 dates (row 2) ---
$ws.Range("Q2").Value = 44221
$ws.Range("R2").Value = 44222
$ws.Range("S2").Value = 44223
$ws.Range("T2").Value = 44224
$ws.Range("U2").Value = 44225
$ws.Range("Q2:U2").NumberFormat = "d-mmm"

# --- Normalize the old "p"/"f" entries in column P to "P"/"F" ---
$colP = @{
    3 = "P"; 4 = "P"; 5 = "P"; 6 = "P"; 7 = "F"; 8 = "F"; 9 = "F"; 10 = "P";
    11 = "F"; 12 = "P"; 13 = "P"; 14 = "P"; 15 = "P"; 16 = "P"; 17 = "P";
    18 = "P"; 19 = "P"; 20 = "P"
}
foreach ($r in 3..20) {
    $ws.Range("P$r").Value = $colP[$r]
}

# --- New attendance column Q for the week-4 block (rows 3-20) ---
$colQ = @{
    3 = "F"; 4 = "F"; 5 = "P"; 6 = "P"; 7 = "F"; 8 = "P"; 9 = "P"; 10 = "P";
    11 = "F"; 12 = "P"; 13 = "F"; 14 = "P"; 15 = "P"; 16 = "P"; 17 = "P";
    18 = "P"; 19 = "P"; 20 = "P"
}
foreach ($r in 3..20) {
    $ws.Range("Q$r").Value = $colQ[$r]
}

# --- Match column width/format of the earlier week blocks ---
$ws.Range("Q1:U20").ColumnWidth = 5.67

# --- Freeze header row + name column, with the active pane on K2 ---
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K2").Select() | Out-Null
